$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 2.9674
$ws.Range("AA2").Value = 2.51345178
$ws.Range("CC2").Value = 55
$ws.Range("CD2").Value = 35
$ws.Range("K3").Value = 7.9026
$ws.Range("K4").Value = 3.3068
$ws.Range("AA4").Value = 0.91515219
$ws.Range("AF4").Value = 77.07317073170699
$ws.Range("AG4").Value = 205
$ws.Range("BG4").Value = 96
$ws.Range("BT4").Value = 47
$ws.Range("DG4").Value = 127
$ws.Range("K5").Value = 2.7085
$ws.Range("AB5").Value = 366
$ws.Range("BB5").Value = 36.781609195402
$ws.Range("BD5").Value = 33.333333333333
$ws.Range("BG5").Value = 121
$ws.Range("CB5").Value = 11
$ws.Range("CU5").Value = 55
$ws.Range("AA6").Value = 1.28663682
$ws.Range("AH6").Value = 248
$ws.Range("AI6").Value = 302
$ws.Range("DF6").Value = 273
$ws.Range("DG6").Value = 363
$ws.Range("K7").Value = 2.6759
$ws.Range("AA7").Value = 0.52050959
$ws.Range("AB7").Value = 294
$ws.Range("AY7").Value = 7
$ws.Range("AZ7").Value = 46.666666666667
$ws.Range("BA7").Value = 47
$ws.Range("BB7").Value = 43.92523364486
$ws.Range("BC7").Value = 46
$ws.Range("BD7").Value = 46.938775510204
$ws.Range("CS7").Value = 15
$ws.Range("AA8").Value = 3.2078558
$ws.Range("AA10").Value = 0.0333936
$ws.Range("CC10").Value = 1
$ws.Range("AA11").Value = 0.2648505
$ws.Range("AW11").Value = 1
$ws.Range("K14").Value = 1.3205
$ws.Range("AA14").Value = 1.05067907
$ws.Range("AB14").Value = 1315
$ws.Range("AF14").Value = 86.00867678958799
$ws.Range("AG14").Value = 922
$ws.Range("AN14").Value = 11.764705882353
$ws.Range("AO14").Value = 16
$ws.Range("AT14").Value = 43
$ws.Range("BG14").Value = 204
$ws.Range("BT14").Value = 129
$ws.Range("CT14").Value = 17
$ws.Range("DF14").Value = 449
$ws.Range("AA15").Value = 1.15789672
$ws.Range("AB15").Value = 735
$ws.Range("AE15").Value = 521
$ws.Range("AF15").Value = 88.155668358714
$ws.Range("AG15").Value = 591
$ws.Range("AH15").Value = 230
$ws.Range("AP15").Value = 25
$ws.Range("AS15").Value = 12
$ws.Range("BA15").Value = 46
$ws.Range("BB15").Value = 47.422680412371
$ws.Range("BC15").Value = 37
$ws.Range("BD15").Value = 47.435897435897
$ws.Range("CU15").Value = 51
$ws.Range("DA15").Value = 16
$ws.Range("DB15").Value = 64
$ws.Range("DF15").Value = 249
$ws.Range("DG15").Value = 342
$ws.Range("K16").Value = 1.874
$ws.Range("AA16").Value = 1.307712
$ws.Range("AF16").Value = 90.104849279161
$ws.Range("AG16").Value = 1526
$ws.Range("AL16").Value = 54.140127388535
$ws.Range("AT16").Value = 116
$ws.Range("BG16").Value = 195
$ws.Range("BT16").Value = 151
$ws.Range("CZ16").Value = 157
$ws.Range("DG16").Value = 668
$ws.Range("AA18").Value = 0.10654326
$ws.Range("AB18").Value = 1074
$ws.Range("AF18").Value = 95.879732739421
$ws.Range("AG18").Value = 898
$ws.Range("AT18").Value = 92
$ws.Range("BA18").Value = 54
$ws.Range("BB18").Value = 56.25
$ws.Range("BE18").Value = 27
$ws.Range("BF18").Value = 56.25
$ws.Range("BG18").Value = 45
$ws.Range("BT18").Value = 37
$ws.Range("DF18").Value = 677
$ws.Range("AA19").Value = 0.53790986
$ws.Range("DA19").Value = 18
$ws.Range("DB19").Value = 69.230769230769
$ws.Range("AA22").Value = 2.14034944
$ws.Range("AE22").Value = 247
$ws.Range("AF22").Value = 74.174174174174
$ws.Range("AG22").Value = 333
$ws.Range("AH22").Value = 89
$ws.Range("AT22").Value = 23
$ws.Range("DF22").Value = 121
$ws.Range("AA26").Value = 0.03342434
$ws.Range("AB26").Value = 817
$ws.Range("BG26").Value = 179
$ws.Range("DL26").Value = -0.5861
